$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 13 de Septiembre de 2020 a las 01:53"

# Row 4: Estados Unidos
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 6674488
$ws.Range("C4").Value = 37169
$ws.Range("D4").Value = 3945557
$ws.Range("E4").Value = 2530813
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 697
$ws.Range("H4").Value = 198118

# Row 6: Brasil
$ws.Range("A6").Value = "Brasil"
$ws.Range("B6").Value = 4315858
$ws.Range("C6").Value = 31880
$ws.Range("D6").Value = 3553421
$ws.Range("E6").Value = 631163
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 800
$ws.Range("H6").Value = 131274

# Row 8: Peru
$ws.Range("A8").Value = "Peru"
$ws.Range("B8").Value = 722832
$ws.Range("C8").Value = 6162
$ws.Range("D8").Value = 559321
$ws.Range("E8").Value = 132918
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 123
$ws.Range("H8").Value = 30593

# Row 24: Alemania
$ws.Range("A24").Value = "Alemania"
$ws.Range("B24").Value = 260546
$ws.Range("C24").Value = 821
$ws.Range("D24").Value = 235300
$ws.Range("E24").Value = 15819
$ws.Range("F24").Value = 0
$ws.Range("G24").Value = 4
$ws.Range("H24").Value = 9427

# Row 36: Panama
$ws.Range("A36").Value = "Panama"
$ws.Range("B36").Value = 101041
$ws.Range("C36").Value = 711
$ws.Range("D36").Value = 73476
$ws.Range("E36").Value = 25410
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = 15
$ws.Range("H36").Value = 2155

# Row 37: Egipto
$ws.Range("A37").Value = "Egipto"
$ws.Range("B37").Value = 100856
$ws.Range("C37").Value = 148
$ws.Range("D37").Value = 83261
$ws.Range("E37").Value = 11968
$ws.Range("F37").Value = 0
$ws.Range("G37").Value = 20
$ws.Range("H37").Value = 5627

# Row 53: Venezuela
$ws.Range("A53").Value = "Venezuela"
$ws.Range("B53").Value = 59630
$ws.Range("C53").Value = 967
$ws.Range("D53").Value = 47729
$ws.Range("E53").Value = 11424
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 9
$ws.Range("H53").Value = 477

# Row 54: Barein
$ws.Range("A54").Value = "Barein"
$ws.Range("B54").Value = 59586
$ws.Range("C54").Value = 747
$ws.Range("D54").Value = 53192
$ws.Range("E54").Value = 6183
$ws.Range("F54").Value = 0
$ws.Range("G54").Value = 3
$ws.Range("H54").Value = 211

# Row 69: Chequia
$ws.Range("A69").Value = "Chequia"
$ws.Range("B69").Value = 35401
$ws.Range("C69").Value = 1541
$ws.Range("D69").Value = 21205
$ws.Range("E69").Value = 13743
$ws.Range("F69").Value = 0
$ws.Range("G69").Value = 3
$ws.Range("H69").Value = 453

# Row 109: Montenegro
$ws.Range("A109").Value = "Montenegro"
$ws.Range("B109").Value = 6530
$ws.Range("C109").Value = 145
$ws.Range("D109").Value = 4491
$ws.Range("E109").Value = 1921
$ws.Range("F109").Value = 0
$ws.Range("G109").Value = 3
$ws.Range("H109").Value = 118

# Row 123: Surinam
$ws.Range("A123").Value = "Surinam"
$ws.Range("B123").Value = 4579
$ws.Range("C123").Value = 50
$ws.Range("D123").Value = 3788
$ws.Range("E123").Value = 698
$ws.Range("F123").Value = 0
$ws.Range("G123").Value = 0
$ws.Range("H123").Value = 93

# Row 124: Ruanda
$ws.Range("A124").Value = "Ruanda"
$ws.Range("B124").Value = 4565
$ws.Range("C124").Value = 31
$ws.Range("D124").Value = 2544
$ws.Range("E124").Value = 1999
$ws.Range("F124").Value = 0
$ws.Range("G124").Value = 0
$ws.Range("H124").Value = 22

# Row 125: Jamaica
$ws.Range("A125").Value = "Jamaica"
$ws.Range("B125").Value = 3623
$ws.Range("C125").Value = 112
$ws.Range("D125").Value = 1072
$ws.Range("E125").Value = 2511
$ws.Range("F125").Value = 0
$ws.Range("G125").Value = 0
$ws.Range("H125").Value = 40

# Row 126: Eslovenia
$ws.Range("A126").Value = "Eslovenia"
$ws.Range("B126").Value = 3603
$ws.Range("C126").Value = 105
$ws.Range("D126").Value = 2699
$ws.Range("E126").Value = 769
$ws.Range("F126").Value = 0
$ws.Range("G126").Value = 0
$ws.Range("H126").Value = 135

# Row 138: Trinidad yTobago
$ws.Range("A138").Value = "Trinidad yTobago"
$ws.Range("B138").Value = 2993
$ws.Range("C138").Value = 168
$ws.Range("D138").Value = 766
$ws.Range("E138").Value = 2176
$ws.Range("F138").Value = 0
$ws.Range("G138").Value = 1
$ws.Range("H138").Value = 51

# Row 139: Bahamas
$ws.Range("A139").Value = "Bahamas"
$ws.Range("B139").Value = 2928
$ws.Range("C139").Value = 114
$ws.Range("D139").Value = 1319
$ws.Range("E139").Value = 1542
$ws.Range("F139").Value = 0
$ws.Range("G139").Value = 2
$ws.Range("H139").Value = 67

# Row 140: Mali
$ws.Range("A140").Value = "Mali"
$ws.Range("B140").Value = 2916
$ws.Range("C140").Value = 4
$ws.Range("D140").Value = 2276
$ws.Range("E140").Value = 512
$ws.Range("F140").Value = 0
$ws.Range("G140").Value = 0
$ws.Range("H140").Value = 128

# Row 149: Islandia
$ws.Range("A149").Value = "Islandia"
$ws.Range("B149").Value = 2162
$ws.Range("C149").Value = 1
$ws.Range("D149").Value = 2085
$ws.Range("E149").Value = 67
$ws.Range("F149").Value = 0
$ws.Range("G149").Value = 0
$ws.Range("H149").Value = 10

# Row 152: Yemen
$ws.Range("A152").Value = "Yemen"
$ws.Range("B152").Value = 2009
$ws.Range("C152").Value = 2
$ws.Range("D152").Value = 1211
$ws.Range("E152").Value = 216
$ws.Range("F152").Value = 0
$ws.Range("G152").Value = 0
$ws.Range("H152").Value = 582

# Row 157: Republica de Chipre
$ws.Range("A157").Value = "Republica de Chipre"
$ws.Range("B157").Value = 1523
$ws.Range("C157").Value = 3
$ws.Range("D157").Value = 1281
$ws.Range("E157").Value = 220
$ws.Range("F157").Value = 0
$ws.Range("G157").Value = 0
$ws.Range("H157").Value = 22

# Row 158: Burkina Faso
$ws.Range("A158").Value = "Burkina Faso"
$ws.Range("B158").Value = 1514
$ws.Range("C158").Value = 15
$ws.Range("D158").Value = 1127
$ws.Range("E158").Value = 331
$ws.Range("F158").Value = 0
$ws.Range("G158").Value = 0
$ws.Range("H158").Value = 56

# Row 191: Monaco
$ws.Range("A191").Value = "Monaco"
$ws.Range("B191").Value = 169
$ws.Range("C191").Value = 1
$ws.Range("D191").Value = 123
$ws.Range("E191").Value = 45
$ws.Range("F191").Value = 0
$ws.Range("G191").Value = 0
$ws.Range("H191").Value = 1

# Row 217: San Pedro y Miquelon
$ws.Range("A217").Value = "San Pedro y Miquelon"
$ws.Range("B217").Value = 11
$ws.Range("C217").Value = 1
$ws.Range("D217").Value = 5
$ws.Range("E217").Value = 6
$ws.Range("F217").Value = 0
$ws.Range("G217").Value = 0
$ws.Range("H217").Value = 0
